$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 1128, shifting existing data (old 1128-1188) down to 1134-1194
$ws.Range("A1128:R1133").EntireRow.Insert()

# Populate the 6 newly inserted rows with the new weekly Cebolla price-report data
# Row 1128
$ws.Range("A1128").Value = 9
$ws.Range("B1128").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1128").Value = "Metropolitana"
$ws.Range("D1128").Value = 44516
$ws.Range("E1128").Value = 13
$ws.Range("F1128").Value = 100112004
$ws.Range("G1128").Value = "Cebolla"
$ws.Range("H1128").Value = "Sin especificar"
$ws.Range("I1128").Value = "1a (guarda)"
$ws.Range("J1128").Value = 340
$ws.Range("K1128").Value = 4600
$ws.Range("L1128").Value = 5000
$ws.Range("M1128").Value = 4800
$ws.Range("N1128").Value = "$/malla 16 kilos"
$ws.Range("O1128").Value = "Región de O'Higgins"
$ws.Range("P1128").Value = 300
$ws.Range("Q1128").Value = 16
$ws.Range("R1128").Value = "Hortaliza"

# Row 1129
$ws.Range("A1129").Value = 9
$ws.Range("B1129").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1129").Value = "Metropolitana"
$ws.Range("D1129").Value = 44516
$ws.Range("E1129").Value = 13
$ws.Range("F1129").Value = 100112004
$ws.Range("G1129").Value = "Cebolla"
$ws.Range("H1129").Value = "Sin especificar"
$ws.Range("I1129").Value = "1a nueva(o)"
$ws.Range("J1129").Value = 340
$ws.Range("K1129").Value = 3800
$ws.Range("L1129").Value = 4200
$ws.Range("M1129").Value = 4000
$ws.Range("N1129").Value = "$/malla 18 kilos"
$ws.Range("O1129").Value = "Región de O'Higgins"
$ws.Range("P1129").Value = 222
$ws.Range("Q1129").Value = 18
$ws.Range("R1129").Value = "Hortaliza"

# Row 1130
$ws.Range("A1130").Value = 9
$ws.Range("B1130").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1130").Value = "Metropolitana"
$ws.Range("D1130").Value = 44516
$ws.Range("E1130").Value = 13
$ws.Range("F1130").Value = 100112004
$ws.Range("G1130").Value = "Cebolla"
$ws.Range("H1130").Value = "Sin especificar"
$ws.Range("I1130").Value = "1a nueva(o)"
$ws.Range("J1130").Value = 8600
$ws.Range("K1130").Value = 1800
$ws.Range("L1130").Value = 2200
$ws.Range("M1130").Value = 2000
$ws.Range("N1130").Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O1130").Value = "Región de O'Higgins"
$ws.Range("P1130").Value = 100
$ws.Range("Q1130").Value = 20
$ws.Range("R1130").Value = "Hortaliza"

# Row 1131
$ws.Range("A1131").Value = 9
$ws.Range("B1131").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1131").Value = "Metropolitana"
$ws.Range("D1131").Value = 44516
$ws.Range("E1131").Value = 13
$ws.Range("F1131").Value = 100112004
$ws.Range("G1131").Value = "Cebolla"
$ws.Range("H1131").Value = "Sin especificar"
$ws.Range("I1131").Value = "2a (guarda)"
$ws.Range("J1131").Value = 196
$ws.Range("K1131").Value = 4000
$ws.Range("L1131").Value = 4400
$ws.Range("M1131").Value = 4200
$ws.Range("N1131").Value = "$/malla 16 kilos"
$ws.Range("O1131").Value = "Región de O'Higgins"
$ws.Range("P1131").Value = 262
$ws.Range("Q1131").Value = 16
$ws.Range("R1131").Value = "Hortaliza"

# Row 1132
$ws.Range("A1132").Value = 9
$ws.Range("B1132").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1132").Value = "Metropolitana"
$ws.Range("D1132").Value = 44516
$ws.Range("E1132").Value = 13
$ws.Range("F1132").Value = 100112004
$ws.Range("G1132").Value = "Cebolla"
$ws.Range("H1132").Value = "Sin especificar"
$ws.Range("I1132").Value = "2a nueva(o)"
$ws.Range("J1132").Value = 142
$ws.Range("K1132").Value = 3200
$ws.Range("L1132").Value = 3600
$ws.Range("M1132").Value = 3400
$ws.Range("N1132").Value = "$/malla 18 kilos"
$ws.Range("O1132").Value = "Región de O'Higgins"
$ws.Range("P1132").Value = 189
$ws.Range("Q1132").Value = 18
$ws.Range("R1132").Value = "Hortaliza"

# Row 1133
$ws.Range("A1133").Value = 9
$ws.Range("B1133").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1133").Value = "Metropolitana"
$ws.Range("D1133").Value = 44516
$ws.Range("E1133").Value = 13
$ws.Range("F1133").Value = 100112004
$ws.Range("G1133").Value = "Cebolla"
$ws.Range("H1133").Value = "Sin especificar"
$ws.Range("I1133").Value = "2a nueva(o)"
$ws.Range("J1133").Value = 3200
$ws.Range("K1133").Value = 1400
$ws.Range("L1133").Value = 1600
$ws.Range("M1133").Value = 1500
$ws.Range("N1133").Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O1133").Value = "Región de O'Higgins"
$ws.Range("P1133").Value = 75
$ws.Range("Q1133").Value = 20
$ws.Range("R1133").Value = "Hortaliza"

